$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (Day serial date, hourly prices, daily avg, slot info)
$ws.Range("A2").Value = 46070
$ws.Range("B2").Value = 4.09
$ws.Range("C2").Value = 2.61
$ws.Range("D2").Value = 0.98
$ws.Range("E2").Value = 0.17
$ws.Range("F2").Value = 0.12
$ws.Range("G2").Value = 1.07
$ws.Range("H2").Value = 1.19
$ws.Range("I2").Value = 1.26
$ws.Range("J2").Value = 4.06
$ws.Range("K2").Value = 8.76
$ws.Range("L2").Value = 1.74
$ws.Range("M2").Value = 0.18
$ws.Range("N2").Value = 0.13
$ws.Range("O2").Value = 0.12
$ws.Range("P2").Value = 0.1
$ws.Range("Q2").Value = 0.1
$ws.Range("R2").Value = 0.13
$ws.Range("S2").Value = 7.83
$ws.Range("T2").Value = 13.11
$ws.Range("U2").Value = 29.04
$ws.Range("V2").Value = 53.03
$ws.Range("W2").Value = 42.43
$ws.Range("X2").Value = 22.98
$ws.Range("Y2").Value = 14.18
$ws.Range("Z2").Value = 8.73

# Slot_4h_max / Slot_4h_price
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 33.16

# Slot_2h_frist / Slot_2h_frist_price
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 47.73

# Slot_2h_second / Slot_2h_second_price
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 21.08

# Slot_min_price
$ws.Range("AG2").Value = "0h-17h"
